$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Append two new rows (5 and 6) that mirror the layout of row 4, adding a new
# "reclamation" record whose "Informations Articles" block spans two lines
# (row5 = first article line, row6 = second "Probleme" line: erreurCommande).
# ---------------------------------------------------------------------------

# Values that look like numbers/dates need to be forced to text so Excel
# doesn't silently convert them (e.g. "0635515554" -> 635515554, or
# "12/08/2024" -> a date serial). Set NumberFormat to "@" (text) on those
# specific cells first, assign the value, then restore formatting from the
# matching cell in row 4 so the final style index matches exactly.
$textCells = @("A5", "C5", "D5", "F5", "J5", "L5")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("A5").Value = "5774457899"
$ws.Range("B5").Value = "Alger"
$ws.Range("C5").Value = "12/08/2024"
$ws.Range("D5").Value = "01/08/2024"
$ws.Range("E5").Value = "LAIB  HAMID"
$ws.Range("F5").Value = "0635515554"
$ws.Range("G5").Value = "erreurLivraison"
$ws.Range("H5").Value = "OPT001CH"
$ws.Range("I5").Value = "INTERRUPTEUR SIMPLE ALLUMAGE"
$ws.Range("J5").Value = "56"
$ws.Range("K5").Value = ""
$ws.Range("L5").Value = "100"
$ws.Range("M5").Value = "Non Conforme"

$ws.Range("G6").Value = "erreurCommande"
$ws.Range("H6").Value = ""
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = ""
$ws.Range("L6").Value = ""
$ws.Range("M6").Value = ""

# Re-apply the exact formatting (border/alignment/style) of row 4 onto the
# new rows 5 and 6, restoring the style indexes that the text-forcing step
# above temporarily changed (copy formats only, values are left untouched).
$ws.Range("A4:F4").Copy()
$ws.Range("A5:F6").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("G4:M4").Copy()
$ws.Range("G5:M6").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("N4").Copy()
$ws.Range("N5:N6").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

# Merge the cells that stay common across the two lines of the new record,
# matching the A1:A3-style merges already used for the header block.
$ws.Range("A5:A6").Merge()
$ws.Range("B5:B6").Merge()
$ws.Range("C5:C6").Merge()
$ws.Range("D5:D6").Merge()
$ws.Range("E5:E6").Merge()
$ws.Range("F5:F6").Merge()
$ws.Range("N5:N6").Merge()
